$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("Z20:AA20").FormulaArray = "={1,2}"
$ws.Range("Z21").FormulaArray = "=_xlfn.ANCHORARRAY(Z20)"
